# Update the BRAE worksheet to split the "Prerequisites" column into
# separate Corequisites / Concurrent / Recommended columns, and move
# "Terms Typically Offered" from column D to column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

# --- A handful of rows need their Prerequisites (column C) text split
# between Prerequisites / Corequisite / Recommended, so the raw text
# changes here too. ----------------------------------------------------
$cUpdates = @(
    @(11, "MATH 142; for engineering students only."),
    @(45, "MATH 242 or MATH 244."),
    @(58, "one of the BIO 111, BIO 161, BOT 121, BRAE 213, or MCRO 221; and CHEM 125 or CHEM 128."),
    @(63, "one of the PHYS 104; PHYS 118; PHYS 121; or PHYS 141; and junior standing.")
)

foreach ($item in $cUpdates) {
    $row = $item[0]
    $text = $item[1]
    $ws.Cells.Item($row, 3).Value = $text
}

# --- Populate the new Corequisites (D) / Concurrent (E) / Recommended
# (F) columns for every data row, and move the old "Terms Typically
# Offered" value (previously column D) into the new column G. ---------
$rowData = @(
    @(2, "NA", "NA", "NA", "F, W"),
    @(3, "NA", "NA", "NA", "F"),
    @(4, "NA", "NA", "NA", "F, W, SP"),
    @(5, "NA", "NA", "NA", "F, SP"),
    @(6, "NA", "NA", "NA", "F"),
    @(7, "NA", "NA", "NA", "F, W"),
    @(8, "NA", "NA", "NA", "W, SP"),
    @(9, "NA", "NA", "NA", "F, W, SP"),
    @(10, "NA", "NA", "NA", "W"),
    @(11, "BIO 213.", "NA", "CHEM 124.", "F, W, SP  "),
    @(12, "NA", "NA", "NA", "F"),
    @(13, "NA", "NA", "NA", "SP"),
    @(14, "NA", "NA", "NA", "SP"),
    @(15, "NA", "NA", "NA", "F"),
    @(16, "NA", "NA", "NA", "F, W, SP"),
    @(17, "NA", "NA", "NA", "F, W, SP"),
    @(18, "NA", "NA", "NA", "F, W, SP"),
    @(19, "NA", "NA", "NA", "W"),
    @(20, "NA", "NA", "NA", "TBD"),
    @(21, "NA", "NA", "NA", "TBD"),
    @(22, "NA", "NA", "NA", "SP"),
    @(23, "NA", "NA", "NA", "W"),
    @(24, "NA", "NA", "NA", "F"),
    @(25, "NA", "NA", "NA", "F"),
    @(26, "NA", "NA", "NA", "SP"),
    @(27, "NA", "NA", "NA", "W"),
    @(28, "NA", "NA", "NA", "W"),
    @(29, "NA", "NA", "NA", "W"),
    @(30, "NA", "NA", "NA", "W"),
    @(31, "NA", "NA", "NA", "SP"),
    @(32, "NA", "NA", "NA", "F"),
    @(33, "NA", "NA", "NA", "W"),
    @(34, "NA", "NA", "NA", "W"),
    @(35, "NA", "NA", "NA", "F, W, SP"),
    @(36, "NA", "NA", "NA", "F, W, SP"),
    @(37, "NA", "NA", "NA", "F"),
    @(38, "NA", "NA", "NA", "W"),
    @(39, "NA", "NA", "NA", "SP"),
    @(40, "NA", "NA", "NA", "F, W"),
    @(41, "NA", "NA", "NA", "F, W, SP"),
    @(42, "NA", "NA", "NA", "F, W, Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and two lower-division courses in GE Area D."),
    @(43, "NA", "NA", "NA", "F, SP"),
    @(44, "NA", "NA", "NA", "F, W, SP"),
    @(45, "STAT 312.", "NA", "NA", "SP "),
    @(46, "NA", "NA", "NA", "SP"),
    @(47, "NA", "NA", "NA", "SP"),
    @(48, "NA", "NA", "NA", "F"),
    @(49, "NA", "NA", "NA", "W"),
    @(50, "NA", "NA", "NA", "F"),
    @(51, "NA", "NA", "NA", "W"),
    @(52, "NA", "NA", "NA", "SP"),
    @(53, "NA", "NA", "NA", "SP"),
    @(54, "NA", "NA", "NA", "SP"),
    @(55, "NA", "NA", "NA", "F"),
    @(56, "NA", "NA", "NA", "W"),
    @(57, "NA", "NA", "NA", "W"),
    @(58, "NA", "NA", "NA", "W"),
    @(59, "NA", "NA", "NA", "W"),
    @(60, "NA", "NA", "NA", "SP"),
    @(61, "NA", "NA", "NA", "SP"),
    @(62, "NA", "NA", "NA", "F"),
    @(63, "NA", "NA", "NA", "SP"),
    @(64, "NA", "NA", "NA", "F"),
    @(65, "NA", "NA", "NA", "SP"),
    @(66, "NA", "NA", "NA", "TBD"),
    @(67, "NA", "NA", "NA", "TBD"),
    @(68, "NA", "NA", "NA", "TBD"),
    @(69, "NA", "NA", "NA", "SP"),
    @(70, "NA", "NA", "NA", "SP"),
    @(71, "NA", "NA", "NA", "TBD"),
    @(72, "NA", "NA", "NA", "F, W, SP"),
    @(73, "NA", "NA", "NA", "SP"),
    @(74, "NA", "NA", "NA", "W"),
    @(75, "NA", "NA", "NA", "TBD"),
    @(76, "NA", "NA", "NA", "TBD"),
    @(77, "NA", "NA", "NA", "F, W, SP")
)

foreach ($item in $rowData) {
    $row = $item[0]
    $coreq = $item[1]
    $concurrent = $item[2]
    $recommended = $item[3]
    $terms = $item[4]

    $ws.Cells.Item($row, 4).Value = $coreq
    $ws.Cells.Item($row, 5).Value = $concurrent
    $ws.Cells.Item($row, 6).Value = $recommended
    $ws.Cells.Item($row, 7).Value = $terms
}
